# Weekly Fruta/Hortalizas update:
# Insert a new weekly observation row right after the current header/first
# data block (row 12), shifting the existing rows 12:24 down to 13:25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 12 (pushes old rows 12-24 down to 13-25)
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with this week's data
$ws.Range("A12").Value = 8
$ws.Range("B12").Value = "Terminal La Palmera de La Serena"
$ws.Range("C12").Value = "Coquimbo"
$ws.Range("D12").Value = 44638
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100101
$ws.Range("H12").Value = "Berries"
$ws.Range("I12").Value = 100101001
$ws.Range("J12").Value = "Arándano (blue)"
$ws.Range("K12").Value = "Sin especificar"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 200
$ws.Range("N12").Value = 6000
$ws.Range("O12").Value = 6500
$ws.Range("P12").Value = 6250
$ws.Range("Q12").Value = "$/bandeja 2 kilos"
$ws.Range("R12").Value = "Provincia de Linares"
$ws.Range("S12").Value = 3125
$ws.Range("T12").Value = 2
